$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "params" (sheet1): fix sckey for Treasury Yield 1Y (t01y)
# ------------------------------------------------------------------
$wsParams = $wb.Worksheets.Item("params")
$wsParams.Cells.Item(64, 6).Value = "DGS1"

# ------------------------------------------------------------------
# Sheet "all-variables" (sheet2)
# ------------------------------------------------------------------
$wsVars = $wb.Worksheets.Item("all-variables")

# Fix sckey for Treasury Yield 1Y (t01y) row as well
$wsVars.Cells.Item(9, 6).Value = "DGS1"

# Mark several rows as initial_forecast = TRUE (column N)
foreach ($r in 3, 19, 20, 21, 24, 26, 27) {
    $wsVars.Cells.Item($r, 14).Value = $true
}

# Expand the table "Table132" to include a new row (row 31)
$lo = $wsVars.ListObjects.Item("Table132")
$lo.Resize($wsVars.Range("A1:Q31"))

# Row 30: fill in "pi" (Personal Income) variable, which already existed as a
# blank placeholder row with some style-only cells
$wsVars.Cells.Item(30, 1).Value = "pi"
$wsVars.Cells.Item(30, 2).Value = "Personal Income"
$wsVars.Cells.Item(30, 3).Value = "Consumer Sector"
$wsVars.Cells.Item(30, 4).Value = "GDP"
$wsVars.Cells.Item(30, 5).Value = "fred"
$wsVars.Cells.Item(30, 6).Value = "RPI"
$wsVars.Cells.Item(30, 7).Value = "bn of 2012 dollars"
$wsVars.Cells.Item(30, 8).Value = "m"
$wsVars.Cells.Item(30, 10).Value = "dlog"
$wsVars.Cells.Item(30, 11).Value = "apchg"
$wsVars.Cells.Item(30, 12).Value = "none"
$wsVars.Cells.Item(30, 14).Value = $true
$wsVars.Cells.Item(30, 15).Value = $false
$wsVars.Cells.Item(30, 16).HorizontalAlignment = -4108

# Copy formatting of row 30 into brand-new row 31 (column by column, skipping
# the "sa" (I) and "append_eom_with_currentval" (M) columns which stay empty)
foreach ($c in 1, 2, 3, 4, 5, 6, 7, 8, 10, 11, 12, 14, 15, 16, 17) {
    $wsVars.Cells.Item(30, $c).Copy()
    $wsVars.Cells.Item(31, $c).PasteSpecial(-4122)
}

# Row 31: new "pid" (Disposable Personal Income) variable
$wsVars.Cells.Item(31, 1).Value = "pid"
$wsVars.Cells.Item(31, 2).Value = "Disposable Personal Income"
$wsVars.Cells.Item(31, 3).Value = "Consumer Sector"
$wsVars.Cells.Item(31, 4).Value = "GDP"
$wsVars.Cells.Item(31, 5).Value = "fred"
$wsVars.Cells.Item(31, 6).Value = "DSPIC96"
$wsVars.Cells.Item(31, 7).Value = "bn of 2012 dollars"
$wsVars.Cells.Item(31, 8).Value = "m"
$wsVars.Cells.Item(31, 10).Value = "dlog"
$wsVars.Cells.Item(31, 11).Value = "apchg"
$wsVars.Cells.Item(31, 12).Value = "none"
$wsVars.Cells.Item(31, 14).Value = $true
$wsVars.Cells.Item(31, 15).Value = $true

# ------------------------------------------------------------------
# View / selection state updates
# ------------------------------------------------------------------

# params: scroll down a bit and select A46:N47
$wsParams.Activate()
$excel.ActiveWindow.ScrollRow = 31
$wsParams.Range("A46:N47").Select()

# all-variables: move selection to N18
$wsVars.Activate()
$wsVars.Range("N18").Select()

# external-forecasts: move selection to A15, and restore it as the active tab
# (matches the workbook's saved activeTab / tabSelected state)
$wsExternal = $wb.Worksheets.Item("external-forecasts")
$wsExternal.Activate()
$wsExternal.Range("A15").Select()

Write-Host "done"
